$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 23:04"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1316878
$ws.Range("C4").Value = 24255
$ws.Range("D4").Value = 221723
$ws.Range("E4").Value = 1016741
$ws.Range("F4").Value = 16723
$ws.Range("G4").Value = 1486
$ws.Range("H4").Value = 78414

# Row 10: Alemania
$ws.Range("A10").Value = "Alemania"
$ws.Range("B10").Value = 170489
$ws.Range("C10").Value = 1059
$ws.Range("D10").Value = 141700
$ws.Range("E10").Value = 21321
$ws.Range("F10").Value = 1712
$ws.Range("G10").Value = 76
$ws.Range("H10").Value = 7468

# Row 16: Peru
$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 61847
$ws.Range("C16").Value = 3321
$ws.Range("D16").Value = 19012
$ws.Range("E16").Value = 41121
$ws.Range("F16").Value = 730
$ws.Range("G16").Value = 87
$ws.Range("H16").Value = 1714

# Row 17: India
$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = 59693
$ws.Range("C17").Value = 3342
$ws.Range("D17").Value = 17883
$ws.Range("E17").Value = 39825
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 96
$ws.Range("H17").Value = 1985

# Row 28: Irlanda
$ws.Range("A28").Value = "Irlanda"
$ws.Range("B28").Value = 22541
$ws.Range("C28").Value = 156
$ws.Range("D28").Value = 17110
$ws.Range("E28").Value = 4002
$ws.Range("F28").Value = 76
$ws.Range("G28").Value = 26
$ws.Range("H28").Value = 1429

# Row 84: Costa de Marfil
$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 1602
$ws.Range("C84").Value = 31
$ws.Range("D84").Value = 754
$ws.Range("E84").Value = 828
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 20

# Row 85: Republica de Macedonia
$ws.Range("A85").Value = "Republica de Macedonia"
$ws.Range("B85").Value = 1586
$ws.Range("C85").Value = 14
$ws.Range("D85").Value = 1099
$ws.Range("E85").Value = 397
$ws.Range("F85").Value = 21
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 90

# Row 105: Niger
$ws.Range("A105").Value = "Niger"
$ws.Range("B105").Value = 795
$ws.Range("C105").Value = 14
$ws.Range("D105").Value = 600
$ws.Range("E105").Value = 151
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 44

# Row 117: Guinea-Bisau
$ws.Range("A117").Value = "Guinea-Bisau"
$ws.Range("B117").Value = 594
$ws.Range("C117").Value = 30
$ws.Range("D117").Value = 25
$ws.Range("E117").Value = 567
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 2

# Row 169: Islas Caimanes
$ws.Range("A169").Value = "Islas Caimanes"
$ws.Range("B169").Value = 81
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 41
$ws.Range("E169").Value = 39
$ws.Range("F169").Value = 3
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1

# Row 205: Seychelles
$ws.Range("A205").Value = "Seychelles"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 8
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Row 206: Montserrat
$ws.Range("A206").Value = "Montserrat"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 7
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 1
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1
